$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-29 Saturday", 2) | Out-Null
$d.Content.Find.Execute("48÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 2) | Out-Null
$d.Content.Find.Execute("66÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("86÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=", 2) | Out-Null
$d.Content.Find.Execute("30÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷5=", 2) | Out-Null
$d.Content.Find.Execute("93÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷6=", 2) | Out-Null
$d.Content.Find.Execute("11÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=", 2) | Out-Null
$d.Content.Find.Execute("11÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷6=", 2) | Out-Null
$d.Content.Find.Execute("84÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷6=", 2) | Out-Null
$d.Content.Find.Execute("71÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=", 2) | Out-Null
$d.Content.Find.Execute("55÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷4=", 2) | Out-Null
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=", 2) | Out-Null
$d.Content.Find.Execute("53÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=", 2) | Out-Null
$d.Content.Find.Execute("13÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=", 2) | Out-Null
$d.Content.Find.Execute("50÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷8=", 2) | Out-Null
$d.Content.Find.Execute("58÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=", 2) | Out-Null
$d.Content.Find.Execute("45÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=", 2) | Out-Null
$d.Content.Find.Execute("44÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=", 2) | Out-Null
$d.Content.Find.Execute("52÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=", 2) | Out-Null
$d.Content.Find.Execute("21÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=", 2) | Out-Null
$d.Content.Find.Execute("54÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2) | Out-Null
$d.Content.Find.Execute("19÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=", 2) | Out-Null
$d.Content.Find.Execute("67÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=", 2) | Out-Null
$d.Content.Find.Execute("82÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=", 2) | Out-Null
$d.Content.Find.Execute("95÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=", 2) | Out-Null
$d.Content.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=", 2) | Out-Null
